$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were added to the "Achicoria" table:
#   - one inserted before the existing row that starts at 2022-04-25 (old row 44)
#   - another inserted three rows later, before the existing row that
#     starts at 2022-06-29 (old row 46, which by then has become row 47)
# Every other existing row simply shifts down to make room (handled by
# EntireRow.Insert below, which preserves the original row contents/styles).

# --- Insert new record #1 at row 44 ---
$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = "2022-10-12"
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112010
$ws.Range("G44").Value = "Achicoria"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 55
$ws.Range("K44").Value = 9000
$ws.Range("L44").Value = 9000
$ws.Range("M44").Value = 9000
$ws.Range("N44").Value = '$/caja 18 unidades'
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 500
$ws.Range("Q44").Value = 18
$ws.Range("R44").Value = "Hortaliza"

# --- Insert new record #2 at row 47 ---
$ws.Rows.Item(47).Insert()

$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = "2022-10-11"
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 100112010
$ws.Range("G47").Value = "Achicoria"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 55
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = 9000
$ws.Range("N47").Value = '$/caja 18 unidades'
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 500
$ws.Range("Q47").Value = 18
$ws.Range("R47").Value = "Hortaliza"
